# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 17:35"

# Row 4 - Madrid
$ws.Range("B4").Value = 67049
$ws.Range("D4").Value = 58118
$ws.Range("E4").Value = 8931

# Row 6 - Castilla y Leon
$ws.Range("B6").Value = 18627
$ws.Range("D6").Value = 16667

# Row 7 - Castilla-La Mancha
$ws.Range("B7").Value = 16789
$ws.Range("D7").Value = 13870
$ws.Range("E7").Value = 2919

# Row 9 - Andalucia
$ws.Range("B9").Value = 12547
$ws.Range("D9").Value = 11172
$ws.Range("E9").Value = 1375

# Row 14 - Ciudad Real
$ws.Range("B14").Value = 5588
$ws.Range("D14").Value = 4740
$ws.Range("E14").Value = 848

# Row 16 - Aragon
$ws.Range("B16").Value = 5195
$ws.Range("D16").Value = 4689
$ws.Range("E16").Value = 506

# Row 20 - La Rioja
$ws.Range("D20").Value = 3679
$ws.Range("E20").Value = 354

# Row 32 - Palencia
$ws.Range("D32").Value = 2067
$ws.Range("E32").Value = 307

# Row 33 - Zamora
$ws.Range("B33").Value = 2307
$ws.Range("D33").Value = 2152
